$d = $word.ActiveDocument

$replacements = @(
    @("933×9=8397", "895×9=8055"),
    @("343×9=3087", "464×7=3248"),
    @("937×9=8433", "695×3=2085"),
    @("368×4=1472", "228×8=1824"),
    @("786×9=7074", "363×5=1815"),
    @("984×4=3936", "627×9=5643"),
    @("638×6=3828", "769×4=3076"),
    @("288×4=1152", "598×9=5382"),
    @("595×2=1190", "840×9=7560"),
    @("238×9=2142", "221×4=884"),
    @("956×4=3824", "341×4=1364"),
    @("299×3=897",  "713×8=5704"),
    @("321×6=1926", "433×3=1299"),
    @("973×7=6811", "413×3=1239"),
    @("917×8=7336", "717×6=4302"),
    @("941×6=5646", "412×8=3296"),
    @("337×7=2359", "211×3=633"),
    @("143×8=1144", "199×4=796"),
    @("451×8=3608", "567×9=5103"),
    @("435×9=3915", "942×5=4710"),
    @("914×6=5484", "754×9=6786"),
    @("359×3=1077", "544×8=4352"),
    @("980×9=8820", "334×4=1336"),
    @("578×6=3468", "847×7=5929"),
    @("523×5=2615", "836×3=2508")
)

foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}
